$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Hamed Daye / Independance Daye) with new title and lyrics
$ws.Range("B2").Value = "Independance Daye"
$ws.Range("C2").Value = 'Jai vu ce quils projettent de faire Ils sont comme des sauterelles Ils déplacent leur civilisation de planète en planète Après avoir consommé toutes les ressources naturelles, ils repartent Nous sommes les suivants Lindépendance daye, lindépendance daye Apocalypse moins trois secondes, cest l''heure H quand sonne le réveil Lindépendance daye, lindépendance daye Cachés sous le soleil de Satan mes soldats veillent Je quitte mon milieu naturel la jungle urbaine quimporte ma dégaine Jerre, jsuis à des années lumières jme sers dmon flair Mon ombre sombre comme Spawn plane sur la ville lumière Dans mon repaire une place pour chaque tête mec, comme trophée de guerre Jpends les boys band aux boyaux des boys scouts qui m''boycottent Check le scoop, branche ton scope découvre les yeux des mioches mortes Les sales gosses en Lacoste prennent la sale dope dans le zen Arrachent ta daube de ma zone ou jte fouette à coup de chaînes Et cache tes pâturages quand vient le prédateur africain Yreste plus rien de la pauvre chèvre et dce vieux monsieur Seguin Jaspire ton atmosphère, pire jentends encore tes soupirs Comme la veille dlapocalypse lève les yeux au ciel, fuir! Mais pour aller où Timal ? viens à la rencontre de clui qui a canné Ltroisième type dun seul coup, mal, mal Mon trip cest dbalourder tes tripes au fond de la trappe J moins rien, les miens se souviendront de cette sale date Lindépendance daye, lindépendance daye Apocalypse moins trois secondes, cest l''heure H quand sonne le réveil Lindépendance daye, lindépendance daye Cachés sous le soleil de Satan mes soldats veillent Jsors dmon cocon, un nouveau colon Sauf que moi mon nom cest pas Christophe Colomb Mais Hamed Daye, le seul chef dcette flotte en perdition Dans ma caravelle jcherche pas dépices Mais le port des MCs qui glissent à fond de cale Pour cette expédition mon drapeau la bannière noire du Secteur Ä flotte Le jeune moussaillon flippe et une flaque Jte fous la pétoche comme mes pères lont foutu jadis aux Boches En foutant une cicatrice sur la joue de la douce France, japproche Alors Papy fait dla résistance et leur mioches collaborent Les pseudos de napalm approchent et veulent ma mort Me foutre en Zone 51, écoute ce baratin Comme un homme-soldat qui a djà échoué, jamais! Comme mon sauce lalbinos, je toast et fausse les données Pour crier nous pas bouger moi et mes têtes grainées Les têtes brûlent encore plus fêlées Qules hommes sous les ordres Pappy Boyington, tché Prêt pour une bataille rangée Lindépendance daye, lindépendance daye Apocalypse moins trois secondes, cest l''heure H quand sonne le réveil Lindépendance daye, lindépendance daye Cachés sous le soleil de Satan mes soldats veillent Pas de mayday, mayday, pas la peine de chialer Jsuis pas le lieutenant Spock, cette loque de captaine Kirk viendra pas te sauver À bord de lEntreprise on nique le bénef Toujours jte fesse, blesse, laisse et jleste, hasta la vista bébé, rhalas Je vois déjà sur ta tête se dresser tes Dread Locks Comme un stressos place toi face à moi et dun seul coup jte boxe Car ya pas de place pour toi dans lhyperespace là où mes forces décuplent On amasse trop de blé, harcelés par les civils en tasse-pé Jme méfie delles et jen foudroie avec mon rayon Delta Appelle moi Seaman, Orco ou bien Chewbacca Quand jenfile mon blouson de combat cest pour la guerre, tché Ton Sud contre not Nord une nouvelle guerre de Sécession Rétorque et sors les armes si tu veux perpétuer la tradition Jviens sauver ton HIP-HOP français comme le Che la fait pour la Bolivie Donc quitte ton baggy, ton treillis ou ton pe-sli, Hamed Daye te suis Lindépendance daye, lindépendance daye Apocalypse moins trois secondes, cest l''heure H quand sonne le réveil Lindépendance daye, lindépendance daye Cachés sous le soleil de Satan mes soldats veillent Encore une fois mon secteur Ä te fout ça Cest chaud daye et tu le sais Hamed Daye, hein, Lino, hein Planque toi Tché, Planque toi Tché Négro on tlavait dit, Secteur Äpocalypse Encore un inédit Lindépendance Daye on tlavait prédit On danse sous une pluie de balles, nique la tendance Le balai tireur Ärsenik pousse le tirailleur, offense Dans toute la France, 98, 99, 2000, 2001, 2002 et tout ce que tu veux !'

# Delete rows 3 through 15, which removes the other songs
$ws.Rows("3:15").Delete()
